$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new bug report row (row 8)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "command use after cryo core activation displays non unique message"
$ws.Range("C8").Value = "ricky"
$ws.Range("D8").Value = "wrong if statement order"
$ws.Range("E8").Value = "unique eventtriggered message added first in if statement"
$ws.Range("F8").Value = "fixed"

# Update selection to new last cell
$ws.Range("F8").Select()

# Update column B width to fit the new (longer) description text
# (runtime adds a fixed padding offset when serializing to XML, so the
# input value is tuned to land exactly on the target stored width of 64)
$ws.Columns.Item(2).ColumnWidth = 63.17
